$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 8.080435
$ws.Range("H2").Value = 24.241305
$ws.Range("I2").Value = 0.1496988574979475
$ws.Range("J2").Value = 0.1496988574979476
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.569028
$ws.Range("N2").Value = 1.707084
$ws.Range("O2").Value = 0.1016535000995941
$ws.Range("P2").Value = 0.1016535000995941
$ws.Range("Q2").Value = 4.597993767179999
$ws.Range("R2").Value = 41.38194390462
$ws.Range("S2").Value = 0.01521741282557673
$ws.Range("T2").Value = 0.01521741282557673
$ws.Range("G3").Value = 8.080435
$ws.Range("H3").Value = 24.241305
$ws.Range("I3").Value = 0.1496988574979475
$ws.Range("J3").Value = 0.1496988574979476
$ws.Range("O3").Value = 0.1962512724671019
$ws.Range("P3").Value = 0.1962512724671019
$ws.Range("Q3").Value = 8.876842673599999
$ws.Range("R3").Value = 79.8915840624
$ws.Range("S3").Value = 0.02937859127084357
$ws.Range("T3").Value = 0.02937859127084357
$ws.Range("G4").Value = 8.080435
$ws.Range("H4").Value = 24.241305
$ws.Range("I4").Value = 0.1496988574979475
$ws.Range("J4").Value = 0.1496988574979476
$ws.Range("M4").Value = 3.083549
$ws.Range("N4").Value = 9.250647000000001
$ws.Range("O4").Value = 0.5508578638987945
$ws.Range("P4").Value = 0.5508578638987945
$ws.Range("Q4").Value = 24.916417263815
$ws.Range("R4").Value = 224.247755374335
$ws.Range("S4").Value = 0.08246279286940943
$ws.Range("T4").Value = 0.08246279286940944
$ws.Range("G5").Value = 8.080435
$ws.Range("H5").Value = 24.241305
$ws.Range("I5").Value = 0.1496988574979475
$ws.Range("J5").Value = 0.1496988574979476
$ws.Range("M5").Value = 0.3400753333333333
$ws.Range("N5").Value = 1.020226
$ws.Range("O5").Value = 0.06075245494223394
$ws.Range("P5").Value = 0.06075245494223393
$ws.Range("Q5").Value = 2.747956626103333
$ws.Range("R5").Value = 24.73160963493
$ws.Range("S5").Value = 0.009094573095047957
$ws.Range("T5").Value = 0.009094573095047957
$ws.Range("G6").Value = 8.080435
$ws.Range("H6").Value = 24.241305
$ws.Range("I6").Value = 0.1496988574979475
$ws.Range("J6").Value = 0.1496988574979476
$ws.Range("M6").Value = 0.5065093333333334
$ws.Range("N6").Value = 1.519528
$ws.Range("O6").Value = 0.0904849085922755
$ws.Range("P6").Value = 0.09048490859227548
$ws.Range("Q6").Value = 4.092815744893334
$ws.Range("R6").Value = 36.83534170404
$ws.Range("S6").Value = 0.01354548743706986
$ws.Range("T6").Value = 0.01354548743706986
$ws.Range("I7").Value = 0.2404784903431001
$ws.Range("J7").Value = 0.2404784903431001
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.569028
$ws.Range("N7").Value = 1.707084
$ws.Range("O7").Value = 0.1016535000995941
$ws.Range("P7").Value = 0.1016535000995941
$ws.Range("Q7").Value = 7.386286162896001
$ws.Range("R7").Value = 66.47657546606401
$ws.Range("S7").Value = 0.02444548024204257
$ws.Range("T7").Value = 0.02444548024204257
$ws.Range("I8").Value = 0.2404784903431001
$ws.Range("J8").Value = 0.2404784903431001
$ws.Range("O8").Value = 0.1962512724671019
$ws.Range("P8").Value = 0.1962512724671019
$ws.Range("Q8").Value = 14.25989323392
$ws.Range("S8").Value = 0.04719420973080109
$ws.Range("T8").Value = 0.04719420973080109
$ws.Range("I9").Value = 0.2404784903431001
$ws.Range("J9").Value = 0.2404784903431001
$ws.Range("M9").Value = 3.083549
$ws.Range("N9").Value = 9.250647000000001
$ws.Range("O9").Value = 0.5508578638987945
$ws.Range("P9").Value = 0.5508578638987945
$ws.Range("Q9").Value = 40.02610646806801
$ws.Range("R9").Value = 360.2349582126121
$ws.Range("S9").Value = 0.132469467504007
$ws.Range("T9").Value = 0.132469467504007
$ws.Range("I10").Value = 0.2404784903431001
$ws.Range("J10").Value = 0.2404784903431001
$ws.Range("M10").Value = 0.3400753333333333
$ws.Range("N10").Value = 1.020226
$ws.Range("O10").Value = 0.06075245494223394
$ws.Range("P10").Value = 0.06075245494223393
$ws.Range("Q10").Value = 4.414358746744001
$ws.Range("R10").Value = 39.72922872069601
$ws.Range("S10").Value = 0.01460965864914563
$ws.Range("T10").Value = 0.01460965864914563
$ws.Range("I11").Value = 0.2404784903431001
$ws.Range("J11").Value = 0.2404784903431001
$ws.Range("M11").Value = 0.5065093333333334
$ws.Range("N11").Value = 1.519528
$ws.Range("O11").Value = 0.0904849085922755
$ws.Range("P11").Value = 0.09048490859227548
$ws.Range("Q11").Value = 6.574760609632001
$ws.Range("R11").Value = 59.172845486688
$ws.Range("S11").Value = 0.02175967421710382
$ws.Range("T11").Value = 0.02175967421710382
$ws.Range("G12").Value = 15.25749233333333
$ws.Range("H12").Value = 45.772477
$ws.Range("I12").Value = 0.2826616599952471
$ws.Range("J12").Value = 0.2826616599952471
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 0.569028
$ws.Range("N12").Value = 1.707084
$ws.Range("O12").Value = 0.1016535000995941
$ws.Range("P12").Value = 0.1016535000995941
$ws.Range("Q12").Value = 8.681940347452
$ws.Range("R12").Value = 78.137463127068
$ws.Range("S12").Value = 0.02873354708247828
$ws.Range("T12").Value = 0.02873354708247828
$ws.Range("G13").Value = 15.25749233333333
$ws.Range("H13").Value = 45.772477
$ws.Range("I13").Value = 0.2826616599952471
$ws.Range("J13").Value = 0.2826616599952471
$ws.Range("O13").Value = 0.1962512724671019
$ws.Range("P13").Value = 0.1962512724671019
$ws.Range("Q13").Value = 16.76127077770667
$ws.Range("R13").Value = 150.85143699936
$ws.Range("S13").Value = 0.05547271045173056
$ws.Range("T13").Value = 0.05547271045173055
$ws.Range("G14").Value = 15.25749233333333
$ws.Range("H14").Value = 45.772477
$ws.Range("I14").Value = 0.2826616599952471
$ws.Range("J14").Value = 0.2826616599952471
$ws.Range("M14").Value = 3.083549
$ws.Range("N14").Value = 9.250647000000001
$ws.Range("O14").Value = 0.5508578638987945
$ws.Range("P14").Value = 0.5508578638987945
$ws.Range("Q14").Value = 47.04722522695767
$ws.Range("R14").Value = 423.425027042619
$ws.Range("S14").Value = 0.1557063982310691
$ws.Range("T14").Value = 0.1557063982310691
$ws.Range("G15").Value = 15.25749233333333
$ws.Range("H15").Value = 45.772477
$ws.Range("I15").Value = 0.2826616599952471
$ws.Range("J15").Value = 0.2826616599952471
$ws.Range("M15").Value = 0.3400753333333333
$ws.Range("N15").Value = 1.020226
$ws.Range("O15").Value = 0.06075245494223394
$ws.Range("P15").Value = 0.06075245494223393
$ws.Range("Q15").Value = 5.188696791089111
$ws.Range("R15").Value = 46.69827111980201
$ws.Range("S15").Value = 0.0171723897627583
$ws.Range("T15").Value = 0.01717238976275829
$ws.Range("G16").Value = 15.25749233333333
$ws.Range("H16").Value = 45.772477
$ws.Range("I16").Value = 0.2826616599952471
$ws.Range("J16").Value = 0.2826616599952471
$ws.Range("M16").Value = 0.5065093333333334
$ws.Range("N16").Value = 1.519528
$ws.Range("O16").Value = 0.0904849085922755
$ws.Range("P16").Value = 0.09048490859227548
$ws.Range("Q16").Value = 7.728062270095112
$ws.Range("R16").Value = 69.552560430856
$ws.Range("S16").Value = 0.02557661446721078
$ws.Range("T16").Value = 0.02557661446721078
$ws.Range("G17").Value = 4.142925
$ws.Range("H17").Value = 12.428775
$ws.Range("I17").Value = 0.07675219702895753
$ws.Range("J17").Value = 0.07675219702895753
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0.569028
$ws.Range("N17").Value = 1.707084
$ws.Range("O17").Value = 0.1016535000995941
$ws.Range("P17").Value = 0.1016535000995941
$ws.Range("Q17").Value = 2.3574403269
$ws.Range("R17").Value = 21.2169629421
$ws.Range("S17").Value = 0.007802129468327199
$ws.Range("T17").Value = 0.007802129468327198
$ws.Range("G18").Value = 4.142925
$ws.Range("H18").Value = 12.428775
$ws.Range("I18").Value = 0.07675219702895753
$ws.Range("J18").Value = 0.07675219702895753
$ws.Range("O18").Value = 0.1962512724671019
$ws.Range("P18").Value = 0.1962512724671019
$ws.Range("Q18").Value = 4.551251688
$ws.Range("R18").Value = 40.961265192
$ws.Range("S18").Value = 0.01506271633157864
$ws.Range("T18").Value = 0.01506271633157863
$ws.Range("G19").Value = 4.142925
$ws.Range("H19").Value = 12.428775
$ws.Range("I19").Value = 0.07675219702895753
$ws.Range("J19").Value = 0.07675219702895753
$ws.Range("M19").Value = 3.083549
$ws.Range("N19").Value = 9.250647000000001
$ws.Range("O19").Value = 0.5508578638987945
$ws.Range("P19").Value = 0.5508578638987945
$ws.Range("Q19").Value = 12.774912240825
$ws.Range("R19").Value = 114.974210167425
$ws.Range("S19").Value = 0.04227955130491095
$ws.Range("T19").Value = 0.04227955130491095
$ws.Range("G20").Value = 4.142925
$ws.Range("H20").Value = 12.428775
$ws.Range("I20").Value = 0.07675219702895753
$ws.Range("J20").Value = 0.07675219702895753
$ws.Range("M20").Value = 0.3400753333333333
$ws.Range("N20").Value = 1.020226
$ws.Range("O20").Value = 0.06075245494223394
$ws.Range("P20").Value = 0.06075245494223393
$ws.Range("Q20").Value = 1.40890660035
$ws.Range("R20").Value = 12.68015940315
$ws.Range("S20").Value = 0.004662884391719204
$ws.Range("T20").Value = 0.004662884391719203
$ws.Range("G21").Value = 4.142925
$ws.Range("H21").Value = 12.428775
$ws.Range("I21").Value = 0.07675219702895753
$ws.Range("J21").Value = 0.07675219702895753
$ws.Range("M21").Value = 0.5065093333333334
$ws.Range("N21").Value = 1.519528
$ws.Range("O21").Value = 0.0904849085922755
$ws.Range("P21").Value = 0.09048490859227548
$ws.Range("Q21").Value = 2.0984301798
$ws.Range("R21").Value = 18.8858716182
$ws.Range("S21").Value = 0.00694491553242154
$ws.Range("T21").Value = 0.006944915532421539
$ws.Range("G22").Value = 13.51654933333334
$ws.Range("H22").Value = 40.549648
$ws.Range("I22").Value = 0.2504087951347477
$ws.Range("J22").Value = 0.2504087951347477
$ws.Range("K22").Value = 3
$ws.Range("L22").Value = 1
$ws.Range("M22").Value = 0.569028
$ws.Range("N22").Value = 1.707084
$ws.Range("O22").Value = 0.1016535000995941
$ws.Range("P22").Value = 0.1016535000995941
$ws.Range("Q22").Value = 7.691295034048001
$ws.Range("R22").Value = 69.22165530643201
$ws.Range("S22").Value = 0.02545493048116931
$ws.Range("T22").Value = 0.02545493048116931
$ws.Range("G23").Value = 13.51654933333334
$ws.Range("H23").Value = 40.549648
$ws.Range("I23").Value = 0.2504087951347477
$ws.Range("J23").Value = 0.2504087951347477
$ws.Range("O23").Value = 0.1962512724671019
$ws.Range("P23").Value = 0.1962512724671019
$ws.Range("Q23").Value = 14.84874043562667
$ws.Range("R23").Value = 133.63866392064
$ws.Range("S23").Value = 0.04914304468214808
$ws.Range("T23").Value = 0.04914304468214807
$ws.Range("G24").Value = 13.51654933333334
$ws.Range("H24").Value = 40.549648
$ws.Range("I24").Value = 0.2504087951347477
$ws.Range("J24").Value = 0.2504087951347477
$ws.Range("M24").Value = 3.083549
$ws.Range("N24").Value = 9.250647000000001
$ws.Range("O24").Value = 0.5508578638987945
$ws.Range("P24").Value = 0.5508578638987945
$ws.Range("Q24").Value = 41.67894218025067
$ws.Range("R24").Value = 375.1104796222561
$ws.Range("S24").Value = 0.1379396539893979
$ws.Range("T24").Value = 0.1379396539893979
$ws.Range("G25").Value = 13.51654933333334
$ws.Range("H25").Value = 40.549648
$ws.Range("I25").Value = 0.2504087951347477
$ws.Range("J25").Value = 0.2504087951347477
$ws.Range("M25").Value = 0.3400753333333333
$ws.Range("N25").Value = 1.020226
$ws.Range("O25").Value = 0.06075245494223394
$ws.Range("P25").Value = 0.06075245494223393
$ws.Range("Q25").Value = 4.596645020049778
$ws.Range("R25").Value = 41.369805180448
$ws.Range("S25").Value = 0.01521294904356285
$ws.Range("T25").Value = 0.01521294904356285
$ws.Range("G26").Value = 13.51654933333334
$ws.Range("H26").Value = 40.549648
$ws.Range("I26").Value = 0.2504087951347477
$ws.Range("J26").Value = 0.2504087951347477
$ws.Range("M26").Value = 0.5065093333333334
$ws.Range("N26").Value = 1.519528
$ws.Range("O26").Value = 0.0904849085922755
$ws.Range("P26").Value = 0.09048490859227548
$ws.Range("Q26").Value = 6.846258391793779
$ws.Range("R26").Value = 61.61632552614401
$ws.Range("S26").Value = 0.02265821693846948
$ws.Range("T26").Value = 0.02265821693846948
